## feat: add 2022-Q1 data
##
## 1) The previous "总计" (totals) sheet is renamed to "2022-Q1" and its
##    contents are replaced with the per-fund holding breakdown for that
##    quarter (same column layout as the other quarterly sheets).
## 2) A brand-new "总计" sheet is appended at the end with the running
##    totals table (one new row for 2022-Q1, prior rows shifted down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while FORCING text storage, even when
# the text looks numeric (e.g. "005585", "0.50") so leading zeros /
# trailing zeros survive instead of Excel auto-coercing to a number.
# Uses a scratch cell formatted as Text, then copies the *value only*
# (PasteSpecial xlPasteValues = -4163) into the destination so the
# destination's own cell style/format is left untouched.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Item(1).Range("ZZ1000")
$scratch.NumberFormat = "@"

function Set-TextValue {
    param($range, $text)
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
}

# -----------------------------------------------------------------
# Step 1: rename "总计" -> "2022-Q1" and overwrite with fund holdings
# -----------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# clear any previous "总计" content before writing the new table
$q1.Cells.Clear()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Rows = @(
    @{ idx = 0; code = "005585"; name = "银河文体娱乐主题灵活配置混合"; scale = "5.54"; pos = "74.07"; ratio = "4.96"; value = "0.2748"; rank = 5 },
    @{ idx = 1; code = "001628"; name = "招商体育文化休闲股票";        scale = "2.95"; pos = "83.21"; ratio = "4.94"; value = "0.1457"; rank = 2 },
    @{ idx = 2; code = "001403"; name = "招商国企改革主题混合";        scale = "3.57"; pos = "87.32"; ratio = "4.03"; value = "0.1439"; rank = 9 },
    @{ idx = 3; code = "002271"; name = "招商安弘灵活配置混合";        scale = "0.50"; pos = "72.34"; ratio = "4.06"; value = "0.0203"; rank = 5 }
)

$row = 2
foreach ($r in $q1Rows) {
    $q1.Range("A$row").Value = $r.idx
    Set-TextValue $q1.Range("B$row") $r.code
    $q1.Range("C$row").Value = $r.name
    Set-TextValue $q1.Range("D$row") $r.scale
    Set-TextValue $q1.Range("E$row") $r.pos
    Set-TextValue $q1.Range("F$row") $r.ratio
    Set-TextValue $q1.Range("G$row") $r.value
    $q1.Range("H$row").Value = $r.rank
    $row = $row + 1
}

# -----------------------------------------------------------------
# Step 2: append a fresh "总计" sheet after the last sheet and fill it
# with the running totals, newest quarter first.
# -----------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @{ idx = 0; date = "2022-Q1"; count = 4; value = 0.58 },
    @{ idx = 1; date = "2021-Q4"; count = 5; value = 1.49 },
    @{ idx = 2; date = "2021-Q2"; count = 1; value = 0.02 }
)

$row = 2
foreach ($r in $totalRows) {
    $total.Range("A$row").Value = $r.idx
    $total.Range("B$row").Value = $r.date
    $total.Range("C$row").Value = $r.count
    $total.Range("D$row").Value = $r.value
    $row = $row + 1
}

# clean up the scratch cell used for forced-text writes
$scratch.Clear()
